# Auto-generated edit script: apply Marilith_Profits scheduled-runner updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 2541.125
$ws.Range("I5").Value = 3378.5
$ws.Range("K5").Value = 3378.5
$ws.Range("M5").Value = -3263.5

$ws.Range("H62").Value = 4284.857
$ws.Range("I62").Value = 4066.5
$ws.Range("J62").Value = 4448.625
$ws.Range("K62").Value = 4066.5
$ws.Range("L62").Value = 4448.625
$ws.Range("M62").Value = -3442.5
$ws.Range("N62").Value = -5696.625

$ws.Range("H65").Value = 4284.857
$ws.Range("I65").Value = 4066.5
$ws.Range("J65").Value = 4448.625
$ws.Range("K65").Value = 20332.5
$ws.Range("L65").Value = 22243.125
$ws.Range("M65").Value = -17212.5
$ws.Range("N65").Value = -28483.125

$ws.Range("H74").Value = 3999
$ws.Range("I74").Value = 3999
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3999
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3063
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 3999
$ws.Range("I77").Value = 3999
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 19995
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -15315
$ws.Range("N77").ClearContents()

$ws.Range("H135").Value = 1359.3846
$ws.Range("I135").Value = 1199.3
$ws.Range("K135").Value = 10793.7
$ws.Range("M135").Value = -8258.699999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 3519
$ws.Range("I26").Value = 3519
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 3519
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -3189
$ws.Range("N26").ClearContents()

$ws.Range("H35").Value = 2871.3333
$ws.Range("I35").Value = 2307.25
$ws.Range("J35").Value = 3999.5
$ws.Range("K35").Value = 2307.25
$ws.Range("L35").Value = 3999.5
$ws.Range("M35").Value = -1901.25
$ws.Range("N35").Value = -4811.5

$ws.Range("H37").Value = 22000
$ws.Range("J37").Value = 25000
$ws.Range("L37").Value = 25000
$ws.Range("N37").Value = -25546

$ws.Range("H45").Value = 1965
$ws.Range("I45").Value = 1965
$ws.Range("K45").Value = 1965
$ws.Range("M45").Value = -1588

$ws.Range("H102").Value = 2203.4
$ws.Range("I102").Value = 2203.4
$ws.Range("K102").Value = 2203.4
$ws.Range("M102").Value = -581.4000000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 101
$ws.Range("J22").Value = 299
$ws.Range("K22").Value = 101
$ws.Range("L22").Value = 299
$ws.Range("M22").Value = 72
$ws.Range("N22").Value = -645

$ws.Range("H95").Value = 26379
$ws.Range("J95").Value = 26379
$ws.Range("L95").Value = 26379
$ws.Range("N95").Value = -31871

$ws.Range("H109").Value = 79986.75
$ws.Range("J109").Value = 79986.75
$ws.Range("L109").Value = 79986.75
$ws.Range("N109").Value = -82760.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 85.92308
$ws.Range("I7").Value = 89.2
$ws.Range("J7").Value = 75
$ws.Range("K7").Value = 89.2
$ws.Range("L7").Value = 75
$ws.Range("M7").Value = 23.8
$ws.Range("N7").Value = -301

$ws.Range("H19").Value = 629.4286
$ws.Range("I19").Value = 389.22223
$ws.Range("J19").Value = 2070.6667
$ws.Range("K19").Value = 389.22223
$ws.Range("L19").Value = 2070.6667
$ws.Range("M19").Value = -219.22223
$ws.Range("N19").Value = -2410.6667

$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 23250.25
$ws.Range("I23").Value = 30333.666
$ws.Range("J23").Value = 2000
$ws.Range("K23").Value = 30333.666
$ws.Range("L23").Value = 2000
$ws.Range("M23").Value = -30093.666
$ws.Range("N23").Value = -2480

$ws.Range("H24").Value = 629.4286
$ws.Range("I24").Value = 389.22223
$ws.Range("J24").Value = 2070.6667
$ws.Range("K24").Value = 389.22223
$ws.Range("L24").Value = 2070.6667
$ws.Range("M24").Value = -219.22223
$ws.Range("N24").Value = -2410.6667

$ws.Range("H26").Value = 3200
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 3200
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 3200
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -3774

$ws.Range("H27").Value = 23250.25
$ws.Range("I27").Value = 30333.666
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 30333.666
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -30141.666
$ws.Range("N27").Value = -2384

$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H81").Value = 34499.75
$ws.Range("J81").Value = 34499.75
$ws.Range("L81").Value = 34499.75
$ws.Range("N81").Value = -36495.75

$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").ClearContents()

$ws.Range("H84").Value = 34499.75
$ws.Range("J84").Value = 34499.75
$ws.Range("L84").Value = 103499.25
$ws.Range("N84").Value = -113483.25

$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").ClearContents()

$ws.Range("H138").Value = 36000
$ws.Range("I138").Value = 36000
$ws.Range("K138").Value = 36000
$ws.Range("M138").Value = -30860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 380.45456
$ws.Range("I2").Value = 149.5
$ws.Range("K2").Value = 897
$ws.Range("M2").Value = -784

$ws.Range("H6").Value = 10204.667
$ws.Range("J6").Value = 29999.5
$ws.Range("L6").Value = 89998.5
$ws.Range("N6").Value = -90224.5

$ws.Range("H33").Value = 49.285713
$ws.Range("J33").Value = 93
$ws.Range("L33").Value = 558
$ws.Range("N33").Value = -1124

$ws.Range("H34").Value = 2257.5789
$ws.Range("I34").Value = 725
$ws.Range("J34").Value = 2666.2666
$ws.Range("K34").Value = 2175
$ws.Range("L34").Value = 7998.7998
$ws.Range("M34").Value = -2091
$ws.Range("N34").Value = -8166.7998

$ws.Range("H62").Value = 7406
$ws.Range("I62").Value = 7406
$ws.Range("K62").Value = 22218
$ws.Range("M62").Value = -21532

$ws.Range("H65").Value = 7406
$ws.Range("I65").Value = 7406
$ws.Range("K65").Value = 66654
$ws.Range("M65").Value = -63222

$ws.Range("H137").Value = 4282.8335
$ws.Range("I137").Value = 3999
$ws.Range("K137").Value = 11997
$ws.Range("M137").Value = -6897

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 561.3333
$ws.Range("I2").Value = 683.3333
$ws.Range("J2").Value = 317.33334
$ws.Range("K2").Value = 683.3333
$ws.Range("L2").Value = 317.33334
$ws.Range("M2").Value = -570.3333
$ws.Range("N2").Value = -543.33334

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1471.7142
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 1583.6666
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 1583.6666
$ws.Range("M22").Value = -505
$ws.Range("N22").Value = -2173.6666

$ws.Range("H27").Value = 1471.7142
$ws.Range("I27").Value = 800
$ws.Range("J27").Value = 1583.6666
$ws.Range("K27").Value = 800
$ws.Range("L27").Value = 1583.6666
$ws.Range("M27").Value = -693
$ws.Range("N27").Value = -1797.6666

$ws.Range("I35").Value = 5000
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 5000
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -4664
$ws.Range("N35").ClearContents()

$ws.Range("H75").Value = 59986.5
$ws.Range("J75").Value = 59986.5
$ws.Range("L75").Value = 59986.5
$ws.Range("N75").Value = -61858.5

$ws.Range("H78").Value = 59986.5
$ws.Range("J78").Value = 59986.5
$ws.Range("L78").Value = 179959.5
$ws.Range("N78").Value = -189319.5

$ws.Range("H94").Value = 22500
$ws.Range("J94").Value = 22500
$ws.Range("L94").Value = 22500
$ws.Range("N94").Value = -23852

$ws.Range("H109").Value = 18598
$ws.Range("J109").Value = 18598
$ws.Range("L109").Value = 18598
$ws.Range("N109").Value = -21372

$ws.Range("H136").Value = 2926.3845
$ws.Range("I136").Value = 2794.8333
$ws.Range("K136").Value = 8384.499899999999
$ws.Range("M136").Value = -5834.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 11499.5
$ws.Range("I45").Value = 8500
$ws.Range("J45").Value = 14499
$ws.Range("K45").Value = 8500
$ws.Range("L45").Value = 14499
$ws.Range("M45").Value = -8009
$ws.Range("N45").Value = -15481

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H101").Value = 18938.75
$ws.Range("J101").Value = 18938.75
$ws.Range("L101").Value = 18938.75
$ws.Range("N101").Value = -25428.75
